$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at N (position 14), mirroring the row
# styles/heights already used by the sheet. This pushes nothing else
# around: the existing M column (which held an "unexpected" extra
# column in this test fixture) stays exactly where it is, and a brand
# new, empty N column opens up next to it.
$ws.Range("N1:N10").EntireColumn.Insert()

# Remember the value that used to live in M3 (the "unexpected" column's
# sole populated cell, "Test") before we repurpose M for national_id.
$oldM3 = $ws.Range("M3").Value2

# New header for column M: national_id
$ws.Range("M1").Value = "national_id"

# Move the old "unexpected" column's data over to the new N column.
$ws.Range("N3").Value = $oldM3

# Populate the national_id value for the data row.
$ws.Range("M3").Value = "ABC123456"
